$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark chapters 21 and 22 (rows 22 and 23, column B) as DONE - reuse the
# exact "DONE" formatting (green fill) already used for the other
# completed chapters by copying the format from an existing DONE cell.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B22").Value = "DONE"
$ws.Range("B23").Value = "DONE"

# Update the active selection to match the saved workbook state
$ws.Range("I24").Select() | Out-Null
